$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 25 and 26 swap coin identity (Polygon <-> WrappedeETH); all other
# rows just get refreshed Price / Volume(1h) figures. Values that look like
# plain numbers are forced back to text (matching the sheets string-typed
# Price/Volume columns) by setting NumberFormat to text before assignment,
# then restoring the default (unstyled) cell style.

$ws.Range('D2').Value = '57.907.30'
$ws.Range('E2').Value = '  -0.56%  '
$ws.Range('D3').Value = '2.449.10'
$ws.Range('E3').Value = '  -2.96%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '524.62'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '131.22'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.40%  '
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.563'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.30%  '
$ws.Range('D9').Value = '2.451.53'
$ws.Range('E9').Value = '  -2.83%  '
$ws.Range('E10').Value = '  -0.29%  '
$ws.Range('E11').Value = '  -2.00%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.96'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -3.90%  '
$ws.Range('E13').Value = '  -2.59%  '
$ws.Range('D14').Value = '2.885.34'
$ws.Range('D15').Value = '57.792.13'
$ws.Range('E15').Value = '  -0.81%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '21.71'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.84%  '
$ws.Range('E17').Value = '  -1.66%  '
$ws.Range('D18').Value = '2.452.02'
$ws.Range('E18').Value = '  -2.96%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.28'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -3.62%  '
$ws.Range('E20').Value = '  -1.39%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '311.05'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -3.47%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.09'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.87%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '64.93'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.41%  '
$ws.Range('B25').Value = 'Polygon'
$ws.Range('C25').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.402'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.23%  '
$ws.Range('B26').Value = 'WrappedeETH'
$ws.Range('C26').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D26').Value = '2.580.59'
$ws.Range('E26').Value = '  -2.00%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.999'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.07%  '
$ws.Range('E28').Value = '  -1.56%  '
$ws.Range('E29').Value = '  -2.29%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '173.37'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.72%  '
$ws.Range('D31').Value = '0.0₃0735'
$ws.Range('E31').Value = '  -2.01%  '
$ws.Range('E32').Value = '  -1.44%  '
$ws.Range('E33').Value = '  -1.39%  '
$ws.Range('E34').Value = '  -4.58%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.998'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.00%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.997'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.10%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '17.81'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.02%  '
$ws.Range('E38').Value = '  -5.11%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.79'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.42%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.818'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +6.03%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '36.21'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.73%  '
$ws.Range('E42').Value = '  -2.40%  '
$ws.Range('E43').Value = '  -1.64%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '261.69'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -5.31%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.584'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.24%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '4.79'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.96%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0918'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.05%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '121.92'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -6.25%  '
$ws.Range('E49').Value = '  -1.18%  '
$ws.Range('E50').Value = '  -1.08%  '
$ws.Range('E51').Value = '  -3.91%  '
